# "add mysql connect info" -- insert 4 new columns (F:I) for the MySQL
# connection fields (SqlPort / SqlName / SqlUser / SqlPwd), repurpose the
# old "Pwd" column (E) into "SqlIP", and fill in the sample row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank columns before the old column F; this shifts the old
# F..J (data-validation placeholder + the col-width-only G:J) to J..N,
# which is exactly what the target workbook shows.
$ws.Range("F1:I1").EntireColumn.Insert()

# --- header row -------------------------------------------------------
$ws.Range("E1").Value = "SqlIP"
$ws.Range("F1").Value = "SqlPort"
$ws.Range("G1").Value = "SqlName"
$ws.Range("H1").Value = "SqlUser"
$ws.Range("I1").Value = "SqlPwd"

# --- sample data row ----------------------------------------------------
$ws.Range("E2").Value = "192.168.0.24"
$ws.Range("F2").Value = 3306
$ws.Range("G2").Value = "app_test"
$ws.Range("H2").Value = "root"
$ws.Range("I2").Value = 123456

# --- column widths (best-fit, recomputed by Excel after the edit) -----
$ws.Columns.Item(1).ColumnWidth = 12
$ws.Columns.Item(4).ColumnWidth = 4.857142857142857
$ws.Columns.Item(5).ColumnWidth = 13.142857142857142
$ws.Columns.Item(6).ColumnWidth = 7.714285714285714
$ws.Columns.Item(7).ColumnWidth = 7.714285714285714
$ws.Columns.Item(8).ColumnWidth = 7.714285714285714
$ws.Columns.Item(9).ColumnWidth = 6.857142857142857

# --- selection ----------------------------------------------------------
$ws.Activate()
$ws.Range("G9").Select()
